$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp title
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 19:52"

# Row 4
$ws.Range("B4").Value = 5672717
$ws.Range("C4").Value = 16743
$ws.Range("D4").Value = 3030902
$ws.Range("E4").Value = 2466163
$ws.Range("G4").Value = 578
$ws.Range("H4").Value = 175652

# Row 6
$ws.Range("B6").Value = 2834755
$ws.Range("C6").Value = 68129
$ws.Range("D6").Value = 2095211
$ws.Range("E6").Value = 685566
$ws.Range("G6").Value = 964
$ws.Range("H6").Value = 53978

# Row 12
$ws.Range("B12").Value = 390037
$ws.Range("C12").Value = 1182
$ws.Range("D12").Value = 364285
$ws.Range("E12").Value = 15174
$ws.Range("G12").Value = 32
$ws.Range("H12").Value = 10578

# Row 13
$ws.Range("B13").Value = 387985
$ws.Range("C13").Value = 3715
$ws.Range("G13").Value = 127
$ws.Range("H13").Value = 28797

# Row 21
$ws.Range("B21").Value = 253108
$ws.Range("C21").Value = 1303
$ws.Range("D21").Value = 233915
$ws.Range("E21").Value = 13154
$ws.Range("G21").Value = 23
$ws.Range("H21").Value = 6039

# Row 27
$ws.Range("B27").Value = 123321
$ws.Range("C27").Value = 167
$ws.Range("D27").Value = 109602
$ws.Range("E27").Value = 4673
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 9046

# Row 35
$ws.Range("B35").Value = 88127
$ws.Range("C35").Value = 1004
$ws.Range("D35").Value = 56760
$ws.Range("E35").Value = 29866
$ws.Range("G35").Value = 12
$ws.Range("H35").Value = 1501

# Row 45
$ws.Range("A45").Value = "Guatemala"
$ws.Range("B45").Value = 64881
$ws.Range("C45").Value = 1034
$ws.Range("D45").Value = 53362
$ws.Range("E45").Value = 9052
$ws.Range("G45").Value = 48
$ws.Range("H45").Value = 2467

# Row 46
$ws.Range("A46").Value = "Paises Bajos"
$ws.Range("B46").Value = 64525
$ws.Range("C46").Value = 552
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("G46").Value = 6
$ws.Range("H46").Value = 6181

# Row 54
$ws.Range("B54").Value = 46313
$ws.Range("C54").Value = 1510
$ws.Range("D54").Value = 31576
$ws.Range("E54").Value = 13994
$ws.Range("G54").Value = 29
$ws.Range("H54").Value = 743

# Row 58
$ws.Range("B58").Value = 39847
$ws.Range("C58").Value = 403
$ws.Range("D58").Value = 27971
$ws.Range("E58").Value = 10474
$ws.Range("G58").Value = 11
$ws.Range("H58").Value = 1402

# Row 70
$ws.Range("B70").Value = 27547
$ws.Range("C70").Value = 48
$ws.Range("E70").Value = 2408

# Row 87
$ws.Range("A87").Value = "Libano"
$ws.Range("B87").Value = 10347
$ws.Range("C87").Value = 589
$ws.Range("D87").Value = 2928
$ws.Range("E87").Value = 7310
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 109

# Row 88
$ws.Range("A88").Value = "Zambia"
$ws.Range("B88").Value = 10218
$ws.Range("C88").Value = 237
$ws.Range("D88").Value = 9126
$ws.Range("E88").Value = 823
$ws.Range("G88").Value = 5
$ws.Range("H88").Value = 269

# Row 89
$ws.Range("A89").Value = "Noruega"
$ws.Range("B89").Value = 10135
$ws.Range("C89").Value = 24
$ws.Range("D89").Value = 8857
$ws.Range("E89").Value = 1016
$ws.Range("H89").Value = 262

# Row 100
$ws.Range("A100").Value = "Grecia"
$ws.Range("B100").Value = 7684
$ws.Range("C100").Value = 212
$ws.Range("D100").Value = 3804
$ws.Range("E100").Value = 3645
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 235

# Row 101
$ws.Range("A101").Value = "Luxemburgo"
$ws.Range("B101").Value = 7499
$ws.Range("D101").Value = 6753
$ws.Range("E101").Value = 622
$ws.Range("H101").Value = 124

# Row 104
$ws.Range("B104").Value = 6225
$ws.Range("C104").Value = 146
$ws.Range("D104").Value = 3788
$ws.Range("E104").Value = 2413

# Row 111
$ws.Range("A111").Value = "Namibia"
$ws.Range("B111").Value = 4665
$ws.Range("C111").Value = 201
$ws.Range("D111").Value = 2426
$ws.Range("E111").Value = 2200
$ws.Range("H111").Value = 39

# Row 112
$ws.Range("A112").Value = "Hong Kong"
$ws.Range("B112").Value = 4587
$ws.Range("C112").Value = 26
$ws.Range("D112").Value = 3779
$ws.Range("E112").Value = 736
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 72

# Row 123
$ws.Range("A123").Value = "Mozambique"
$ws.Range("B123").Value = 3045
$ws.Range("C123").Value = 54
$ws.Range("D123").Value = 1291
$ws.Range("E123").Value = 1735
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 19

# Row 124
$ws.Range("A124").Value = "Eslovaquia"
$ws.Range("B124").Value = 3022
$ws.Range("C124").Value = 100
$ws.Range("D124").Value = 1997
$ws.Range("E124").Value = 992
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 33

# Row 147
$ws.Range("A147").Value = "Republica de Chipre"
$ws.Range("B147").Value = 1385
$ws.Range("C147").Value = 26
$ws.Range("D147").Value = 878
$ws.Range("E147").Value = 487
$ws.Range("H147").Value = 20

# Row 148
$ws.Range("A148").Value = "Georgia"
$ws.Range("B148").Value = 1361
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 1098
$ws.Range("E148").Value = 246
$ws.Range("H148").Value = 17

# Row 151
$ws.Range("A151").Value = "Burkina Faso"
$ws.Range("B151").Value = 1285
$ws.Range("C151").Value = 5
$ws.Range("D151").Value = 1023
$ws.Range("E151").Value = 207
$ws.Range("H151").Value = 55

# Row 152
$ws.Range("A152").Value = "Liberia"
$ws.Range("B152").Value = 1282
$ws.Range("D152").Value = 803
$ws.Range("E152").Value = 397
$ws.Range("H152").Value = 82

# Row 174
$ws.Range("B174").Value = 382
$ws.Range("C174").Value = 5
$ws.Range("D174").Value = 256
$ws.Range("E174").Value = 126

# Row 195
$ws.Range("B195").Value = 94
$ws.Range("C195").Value = 1
$ws.Range("E195").Value = 3

# Row 213
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
